$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 46011
$ws.Range("I11").Value = 46011
$ws.Range("K11").Value = 46011
$ws.Range("M11").Value = -45871
# Row 15
$ws.Range("H15").Value = 2513.3914
$ws.Range("I15").Value = 2513.3914
$ws.Range("K15").Value = 7540.174199999999
$ws.Range("M15").Value = -7371.174199999999
# Row 51
$ws.Range("H51").Value = 2945.6
$ws.Range("I51").Value = 2945.6
$ws.Range("K51").Value = 2945.6
$ws.Range("M51").Value = -2461.6
# Row 74
$ws.Range("H74").Value = 9722.223
$ws.Range("I74").Value = 9722.223
$ws.Range("K74").Value = 9722.223
$ws.Range("M74").Value = -8786.223
# Row 77
$ws.Range("H77").Value = 9722.223
$ws.Range("I77").Value = 9722.223
$ws.Range("K77").Value = 48611.115
$ws.Range("M77").Value = -43931.115
# Row 100
$ws.Range("H100").Value = 2221.2222
$ws.Range("I100").Value = 1998.875
$ws.Range("J100").Value = 4000
$ws.Range("K100").Value = 1998.875
$ws.Range("L100").Value = 4000
$ws.Range("M100").Value = -1457.875
$ws.Range("N100").Value = -5082
# Row 106
$ws.Range("H106").Value = 6511.154
$ws.Range("I106").Value = 5387.1665
$ws.Range("J106").Value = 19999
$ws.Range("K106").Value = 5387.1665
$ws.Range("L106").Value = 19999
$ws.Range("M106").Value = -4756.1665
$ws.Range("N106").Value = -21261
# Row 113
$ws.Range("H113").Value = 5165.5
$ws.Range("J113").Value = 5748.25
$ws.Range("L113").Value = 5748.25
$ws.Range("N113").Value = -12256.25
# Row 116
$ws.Range("H116").Value = 27822.5
$ws.Range("I116").Value = 26845.6
$ws.Range("K116").Value = 26845.6
$ws.Range("M116").Value = -23403.6
# Row 138
$ws.Range("H138").Value = 4366.577
$ws.Range("J138").Value = 4518.091
$ws.Range("L138").Value = 13554.273
$ws.Range("N138").Value = -23834.273

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1105.5
$ws.Range("I2").Value = 958.06665
$ws.Range("J2").Value = 1547.8
$ws.Range("K2").Value = 958.06665
$ws.Range("L2").Value = 1547.8
$ws.Range("M2").Value = -845.06665
$ws.Range("N2").Value = -1773.8
# Row 32
$ws.Range("H32").Value = 960678.5600000001
$ws.Range("I32").Value = 1090422.1
$ws.Range("J32").Value = 15404.429
$ws.Range("K32").Value = 1090422.1
$ws.Range("L32").Value = 15404.429
$ws.Range("M32").Value = -1090135.1
$ws.Range("N32").Value = -15978.429
# Row 45
$ws.Range("H45").Value = 2622.8333
$ws.Range("I45").Value = 2027.4
$ws.Range("K45").Value = 2027.4
$ws.Range("M45").Value = -1650.4
# Row 110
$ws.Range("H110").Value = 2259.0715
$ws.Range("J110").Value = 2000
$ws.Range("L110").Value = 2000
$ws.Range("N110").Value = -6090
# Row 116
$ws.Range("H116").Value = 1105.5
$ws.Range("I116").Value = 958.06665
$ws.Range("J116").Value = 1547.8
$ws.Range("K116").Value = 958.06665
$ws.Range("L116").Value = 1547.8
$ws.Range("M116").Value = 1335.93335
$ws.Range("N116").Value = -6135.8
# Row 132
$ws.Range("H132").Value = 7289.2334
$ws.Range("I132").Value = 5918.0625
$ws.Range("K132").Value = 17754.1875
$ws.Range("M132").Value = -15224.1875
# Row 134
$ws.Range("H134").Value = 68833.336
$ws.Range("J134").Value = 68833.336
$ws.Range("L134").Value = 68833.336
$ws.Range("N134").Value = -78973.336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1105.5
$ws.Range("I3").Value = 958.06665
$ws.Range("J3").Value = 1547.8
$ws.Range("K3").Value = 958.06665
$ws.Range("L3").Value = 1547.8
$ws.Range("M3").Value = -844.06665
$ws.Range("N3").Value = -1775.8
# Row 64
$ws.Range("H64").Value = 35216.668
$ws.Range("I64").Value = 100000
$ws.Range("J64").Value = 2825
$ws.Range("K64").Value = 100000
$ws.Range("L64").Value = 2825
$ws.Range("M64").Value = -99775
$ws.Range("N64").Value = -3275
# Row 67
$ws.Range("H67").Value = 35216.668
$ws.Range("I67").Value = 100000
$ws.Range("J67").Value = 2825
$ws.Range("K67").Value = 100000
$ws.Range("L67").Value = 2825
$ws.Range("M67").Value = -99220
$ws.Range("N67").Value = -4385
# Row 82
$ws.Range("H82").Value = 22164.857
$ws.Range("J82").Value = 46132.668
$ws.Range("L82").Value = 46132.668
$ws.Range("N82").Value = -46898.668
# Row 85
$ws.Range("H85").Value = 22164.857
$ws.Range("J85").Value = 46132.668
$ws.Range("L85").Value = 46132.668
$ws.Range("N85").Value = -48784.668
# Row 99
$ws.Range("H99").Value = 13976.25
$ws.Range("I99").Value = 15758.571
$ws.Range("K99").Value = 15758.571
$ws.Range("M99").Value = -14260.571
# Row 107
$ws.Range("H107").Value = 1611.6945
$ws.Range("I107").Value = 1341.909
$ws.Range("K107").Value = 1341.909
$ws.Range("M107").Value = 578.0909999999999
# Row 134
$ws.Range("H134").Value = 12825512
$ws.Range("I134").Value = 5415.4165
$ws.Range("K134").Value = 16246.2495
$ws.Range("M134").Value = -13711.2495

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 197.28572
$ws.Range("I7").Value = 125
$ws.Range("K7").Value = 125
$ws.Range("M7").Value = -12
# Row 16
$ws.Range("H16").Value = 14712.333
$ws.Range("I16").Value = 18515.857
$ws.Range("K16").Value = 18515.857
$ws.Range("M16").Value = -18228.857
# Row 50
$ws.Range("H50").Value = 34466.668
$ws.Range("J50").Value = 34466.668
$ws.Range("L50").Value = 34466.668
$ws.Range("N50").Value = -35716.668
# Row 107
$ws.Range("H107").Value = 1239.909
$ws.Range("I107").Value = 1061.6666
$ws.Range("J107").Value = 1453.8
$ws.Range("K107").Value = 1061.6666
$ws.Range("L107").Value = 1453.8
$ws.Range("M107").Value = 858.3334
$ws.Range("N107").Value = -5293.8
# Row 112
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = 0
# Row 113
$ws.Range("H113").Value = 14712.333
$ws.Range("I113").Value = 18515.857
$ws.Range("K113").Value = 18515.857
$ws.Range("M113").Value = -16345.857
# Row 118
$ws.Range("H118").Value = 59748.4
$ws.Range("J118").Value = 59748.4
$ws.Range("L118").Value = 59748.4
$ws.Range("N118").Value = -63062.4
# Row 119
$ws.Range("H119").Value = 66387.5
$ws.Range("J119").Value = 66387.5
$ws.Range("L119").Value = 66387.5
$ws.Range("N119").Value = -76063.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Range("H9").Value = 768141.5600000001
$ws.Range("J9").Value = 768141.5600000001
$ws.Range("L9").Value = 2304424.68
$ws.Range("N9").Value = -2304872.68
# Row 37
$ws.Range("H37").Value = 46578.42
$ws.Range("J37").Value = 46578.42
$ws.Range("L37").Value = 139735.26
$ws.Range("N37").Value = -139959.26
# Row 51
$ws.Range("I51").Value = 200002000
$ws.Range("J51").Value = 2996
$ws.Range("K51").Value = 600006000
$ws.Range("L51").Value = 8988
$ws.Range("M51").Value = -600005540
$ws.Range("N51").Value = -9908
# Row 122
$ws.Range("H122").Value = 850543.2
$ws.Range("J122").Value = 1728.4286
$ws.Range("L122").Value = 15555.8574
$ws.Range("N122").Value = -20455.8574

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 14358.137
$ws.Range("I70").Value = 20546.385
$ws.Range("J70").Value = 5419.5557
$ws.Range("K70").Value = 20546.385
$ws.Range("L70").Value = 5419.5557
$ws.Range("M70").Value = -20276.385
$ws.Range("N70").Value = -5959.5557
# Row 73
$ws.Range("H73").Value = 14358.137
$ws.Range("I73").Value = 20546.385
$ws.Range("J73").Value = 5419.5557
$ws.Range("K73").Value = 20546.385
$ws.Range("L73").Value = 5419.5557
$ws.Range("M73").Value = -19610.385
$ws.Range("N73").Value = -7291.5557
# Row 113
$ws.Range("H113").Value = 1964.9412
$ws.Range("J113").Value = 2899.6667
$ws.Range("L113").Value = 2899.6667
$ws.Range("N113").Value = -7239.6667
# Row 122
$ws.Range("H122").Value = 3175.6562
$ws.Range("I122").Value = 3320.9333
$ws.Range("K122").Value = 9962.7999
$ws.Range("M122").Value = -7512.7999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 300
$ws.Range("I22").Value = 300
$ws.Range("K22").Value = 300
$ws.Range("M22").Value = -5
# Row 27
$ws.Range("H27").Value = 300
$ws.Range("I27").Value = 300
$ws.Range("K27").Value = 300
$ws.Range("M27").Value = -193
# Row 68
$ws.Range("H68").Value = 5037.357
$ws.Range("I68").Value = 4007.4167
$ws.Range("J68").Value = 11217
$ws.Range("K68").Value = 4007.4167
$ws.Range("L68").Value = 11217
$ws.Range("M68").Value = -3258.4167
$ws.Range("N68").Value = -12715
# Row 71
$ws.Range("H71").Value = 5037.357
$ws.Range("I71").Value = 4007.4167
$ws.Range("J71").Value = 11217
$ws.Range("K71").Value = 20037.0835
$ws.Range("L71").Value = 56085
$ws.Range("M71").Value = -16293.0835
$ws.Range("N71").Value = -63573

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 136
$ws.Range("H136").Value = 4078506.5
$ws.Range("I136").Value = 1892798.8
$ws.Range("K136").Value = 5678396.4
$ws.Range("M136").Value = -5675846.4
